# Applies three changes described by the diff:
# 1. Insert a new "Meta description: ..." paragraph right after the
#    Heading1 title paragraph (top of document).
# 2. Remove the paragraph that duplicated the title in bold near the
#    bottom of the document ("Play Book of Tombs Free - ...").
# 3. Replace the text of the remaining (italic) paragraph that used to
#    hold the review blurb with the new image-prompt text, keeping the
#    italic run formatting intact.

$d = $word.ActiveDocument

# --- Step 1: insert the "Meta description" paragraph after the title ---

$titlePara = $d.Paragraphs.First
$titleEnd = $titlePara.Range
$titleEnd.Collapse(0)            # wdCollapseEnd
$titleEnd.InsertParagraphAfter() # creates a new (empty) paragraph #2

$metaPara = $d.Paragraphs.Item(2)
$insertRng = $metaPara.Range
$insertRng.Collapse(1)           # wdCollapseStart

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Book of Tombs, an Egyptian-themed slot game. Play for free and enjoy exceptional graphics, high maximum payouts and a free spin bonus.</w:t></w:r></w:p>'
[void]$insertRng.InsertXML($metaXml)

# --- Step 2: delete the duplicated bold "Play Book of Tombs Free..." paragraph ---
# (search from the bottom so the Heading1 title at the very top is never touched;
#  also require the paragraph style to be Normal, i.e. not the real heading)

$boldText = "Play Book of Tombs Free - Egyptian-themed Slot Game"
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq $boldText -and $para.Style.NameLocal -eq "Normal") {
        $para.Range.Delete()
        break
    }
}

# --- Step 3: replace the italic review-blurb paragraph text with the new prompt ---

$oldBlurb = "Read our review of Book of Tombs, an Egyptian-themed slot game. Play for free and enjoy exceptional graphics, high maximum payouts and a free spin bonus."
$newBlurb = 'Create a cartoon-style feature image for the game "Book of Tombs" that features a happy Maya warrior with glasses. The image should be vibrant and eye-catching, using warm colors to evoke the Egyptian theme of the game. The Maya warrior should be shown holding the book of the Pharaoh, with a confident expression on their face. The background should feature the pyramids and other Egyptian landmarks, with the logo of the game prominently displayed. It should convey a sense of adventure and excitement, making players eager to dive into the game and uncover the treasures that await them.'

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq $oldBlurb) {
        $r = $para.Range
        $textRng = $d.Range($r.Start, $r.End - 1)   # exclude paragraph mark
        $textRng.Text = $newBlurb
        break
    }
}

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
